# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F2").Value = 119
$sheet1.Range("F3").Value = 537
$sheet1.Range("F6").Value = 497
$sheet1.Range("F7").Value = 97
$sheet1.Range("F8").Value = 113
$sheet1.Range("F10").Value = 6625
$sheet1.Range("F11").Value = 226
$sheet1.Range("F13").Value = 2845
$sheet1.Range("F14").Value = 176
$sheet1.Range("F15").Value = 306
$sheet1.Range("F16").Value = 257
$sheet1.Range("F17").Value = 525

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F4").Value = 119
$sheet4.Range("F5").Value = 537
$sheet4.Range("F8").Value = 497
$sheet4.Range("F9").Value = 97
$sheet4.Range("F10").Value = 113
$sheet4.Range("F13").Value = 6625
$sheet4.Range("F15").Value = 226
$sheet4.Range("F17").Value = 2845
$sheet4.Range("F18").Value = 176
$sheet4.Range("F19").Value = 306
$sheet4.Range("F20").Value = 257
$sheet4.Range("F21").Value = 525
